# Fixed a bug in WinWeight.Fit
# The data rows (2-25) of Sheet1 need to be reordered/rewritten so that
# each row's symbol+reel values match the corrected fit ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(801,3,67,65,52,45),
    @(1001,18,30,75,60,72),
    @(701,3,90,45,97,15),
    @(201,9,30,15,45,30),
    @(401,9,48,67,75,45),
    @(1203,3,15,15,15,15),
    @(501,9,52,30,75,45),
    @(601,9,60,67,60,42),
    @(1202,2,10,10,10,10),
    @(301,6,45,30,60,45),
    @(1201,2,10,10,10,10),
    @(101,9,30,15,60,15),
    @(901,16,15,45,60,60),
    @(902,1,0,0,0,0),
    @(502,0,4,0,0,0),
    @(1,0,2,2,2,2),
    @(1101,0,15,30,30,0),
    @(802,0,4,5,4,0),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(602,0,0,4,0,9),
    @(402,0,0,4,0,0),
    @(702,0,0,0,4,0),
    @(1002,0,0,0,0,9)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
